$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row containing account 004500804 / RAFAEL / 5000 (Excel row 7),
# shifting all rows below it up by one.
$ws.Rows.Item(7).Delete()
